# Auto-generated Excel COM-interop script to apply odds updates
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("Q2").Value = 1.88
$ws.Range("R2").Value = 1.98
$ws.Range("S2").Value = 2.5
$ws.Range("T2").Value = 1.5

# Row 3
$ws.Range("Q3").Value = 3.2
$ws.Range("R3").Value = 1.36
$ws.Range("U3").Value = 8.4
$ws.Range("V3").Value = 1.07

# Row 5
$ws.Range("M5").Value = 1.11
$ws.Range("N5").Value = 6.5
$ws.Range("U5").Value = 4.5
$ws.Range("V5").Value = 1.21

# Row 6
$ws.Range("G6").Value = 2.45
$ws.Range("I6").Value = 3.5
$ws.Range("AQ6").Value = 41

# Row 8
$ws.Range("M8").Value = 1.1
$ws.Range("N8").Value = 7
$ws.Range("Y8").Value = 1.53
$ws.Range("Z8").Value = 2.38
$ws.Range("AM8").Value = 10
$ws.Range("AN8").Value = 21

# Row 10
$ws.Range("G10").Value = 2.77
$ws.Range("H10").Value = 2.52
$ws.Range("J10").Value = 3.6
$ws.Range("L10").Value = 3.9
$ws.Range("N10").Value = 4.4
$ws.Range("O10").Value = 1.7
$ws.Range("P10").Value = 2.05
$ws.Range("S10").Value = 3
$ws.Range("T10").Value = 1.34
$ws.Range("W10").Value = 5.5
$ws.Range("X10").Value = 1.11
$ws.Range("Y10").Value = 1.7
$ws.Range("Z10").Value = 2.05
$ws.Range("AC10").Value = 5.7
$ws.Range("AD10").Value = 12
$ws.Range("AE10").Value = 11.25
$ws.Range("AG10").Value = 35
$ws.Range("AH10").Value = 60
$ws.Range("AI10").Value = 4.4
$ws.Range("AM10").Value = 6.2
$ws.Range("AO10").Value = 12

# Row 11
$ws.Range("G11").Value = 2.67
$ws.Range("H11").Value = 2.55
$ws.Range("I11").Value = 3.15
$ws.Range("J11").Value = 3.6
$ws.Range("L11").Value = 3.95
$ws.Range("M11").Value = 1.19
$ws.Range("N11").Value = 4.2
$ws.Range("AD11").Value = 11
$ws.Range("AE11").Value = 12
$ws.Range("AG11").Value = 37
$ws.Range("AH11").Value = 75
$ws.Range("AI11").Value = 4.2
$ws.Range("AJ11").Value = 5.5
$ws.Range("AK11").Value = 22
$ws.Range("AM11").Value = 6.2
$ws.Range("AO11").Value = 12.5
$ws.Range("AQ11").Value = 40
$ws.Range("AR11").Value = 70

# Row 12
$ws.Range("G12").Value = 2.88
$ws.Range("H12").Value = 3.5
$ws.Range("I12").Value = 2.25
$ws.Range("K12").Value = 2.3
$ws.Range("L12").Value = 2.88
$ws.Range("AF12").Value = 29
$ws.Range("AK12").Value = 11
$ws.Range("AQ12").Value = 17
$ws.Range("AS12").Value = 101

# Row 17
$ws.Range("G17").Value = 3.3
$ws.Range("I17").Value = 2
$ws.Range("M17").Value = 1.02
$ws.Range("N17").Value = 12
$ws.Range("O17").Value = 1.25
$ws.Range("P17").Value = 3.75
$ws.Range("S17").Value = 1.8
$ws.Range("T17").Value = 2
$ws.Range("AE17").Value = 12
$ws.Range("AF17").Value = 34
$ws.Range("AN17").Value = 11
$ws.Range("AP17").Value = 19
$ws.Range("AQ17").Value = 17

# Row 18
$ws.Range("G18").Value = 2.8
$ws.Range("I18").Value = 2.3
$ws.Range("J18").Value = 3.1
$ws.Range("L18").Value = 2.75
$ws.Range("M18").Value = 1.01
$ws.Range("N18").Value = 15
$ws.Range("AC18").Value = 13
$ws.Range("AD18").Value = 17
$ws.Range("AF18").Value = 29
$ws.Range("AG18").Value = 21
$ws.Range("AM18").Value = 11
$ws.Range("AN18").Value = 13
$ws.Range("AQ18").Value = 17
$ws.Range("AR18").Value = 21

# Row 20
$ws.Range("S20").Value = 2.08
$ws.Range("T20").Value = 1.73
$ws.Range("W20").Value = 3.75
$ws.Range("X20").Value = 1.25

# Row 24
$ws.Range("G24").Value = 1.53
$ws.Range("H24").Value = 4.25
$ws.Range("J24").Value = 2.02
$ws.Range("K24").Value = 2.42
$ws.Range("N24").Value = 9
$ws.Range("O24").Value = 1.18
$ws.Range("P24").Value = 4.35
$ws.Range("S24").Value = 1.53
$ws.Range("T24").Value = 2.32
$ws.Range("W24").Value = 2.3
$ws.Range("X24").Value = 1.55
$ws.Range("Y24").Value = 1.29
$ws.Range("Z24").Value = 3.3
$ws.Range("AA24").Value = 1.65
$ws.Range("AB24").Value = 2.12
$ws.Range("AC24").Value = 9
$ws.Range("AD24").Value = 8.5
$ws.Range("AE24").Value = 8
$ws.Range("AG24").Value = 11.25
$ws.Range("AH24").Value = 20
$ws.Range("AI24").Value = 9
$ws.Range("AJ24").Value = 8.5
$ws.Range("AM24").Value = 18.5
